$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(196).Delete()

$ws.Range("A190").Select()
$ws.Range("C198").Select()
